# Update the "update" points column (C2: 25.1 -> 26.1) and strip its number-format style
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 26.1
$ws.Range("C2").ClearFormats()

# Append the new ranking log rows (295-317)
$names = @(
    "יהלי גודר",
    "איתי בסטקר",
    "תאיו ורד",
    "ליהי בראל",
    "אביב ואסקז",
    "איתי בסטקר",
    "אורי שטרנברג",
    "איתי הראל",
    "תומר ששון",
    "תומר ששון",
    "איתי בסטקר",
    "רומי הרשקוביץ",
    "דן פימה",
    "גלי זליג",
    "איתי הראל",
    "יולי יערי תליו",
    "הילה שולויס",
    "ליאם דיין",
    "יהלי דוייב",
    "ירון גלפנד",
    "אן מרש",
    "יהלי דוייב",
    "גלי זליג"
)

$points = @(1,1,6,6,1,1,1,1,1,6,6,1,1,1,1,1,1,1,1,1,1,6,6)

$startRow = 295
for ($i = 0; $i -lt $names.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $names[$i]
    $ws.Cells.Item($r, 2).Value = $points[$i]
}

# Restore the view state (scrolled position + selection) seen after the edit
$ws.Range("A306").Select() | Out-Null
